$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 9; $r++) {
    $ws.Range("E$r").Value = "vnorlr2327523050"
    $ws.Range("I$r").Value = "'9226624645"
    $ws.Range("AX$r").Value = "'3798883283"
}
